$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.308.75'
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").Value = '2.649.65'
$ws.Range("E3").Value = '  +2.01%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("E9").Value = '  +7.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.401'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.85'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.56%  '
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.11%  '
$ws.Range("D14").Value = '3.127.32'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000179'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +16.90%  '
$ws.Range("D16").Value = '65.151.83'
$ws.Range("E16").Value = '  +2.82%  '
$ws.Range("D17").Value = '2.664.57'
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.82%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.88%  '
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.166'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.84%  '
$ws.Range("D31").Value = '0.0₃0930'
$ws.Range("E31").Value = '  +9.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.430'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.16'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.61%  '
$ws.Range("E44").Value = '  +4.63%  '
$ws.Range("E45").Value = '  +4.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.652'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.57%  '
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("E50").Value = '  +2.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.31%  '

# Rows 29 and 30 swap coins (Bittensor moves up, Binance-PegBSC-USD moves down)
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '540.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.65%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.77%  '
